$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.326164245605469
$ws.Range("B1").Value = 4.610901832580566
$ws.Range("C1").Value = 3.482093811035156
$ws.Range("D1").Value = 2.399358510971069
$ws.Range("E1").Value = 2.145503997802734
